$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "swfwewfd"

$ws.Range("B1").Value = 1691
$ws.Range("E1").Value = 70
$ws.Range("H1").Value = 95.85798816568047
$ws.Range("I1").Value = 0.0437731196054254
$ws.Range("J1").Value = 51.03661918640137
$ws.Range("J2").Value = 36.93914771080017
$ws.Range("J3").Value = 44.37130188941956
$ws.Range("B4").Value = 2942
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 99.9659979598776
$ws.Range("I4").Value = 0.01143241425689307
$ws.Range("J4").Value = 41.87565040588379
$ws.Range("J5").Value = 44.60436964035034
$ws.Range("J6").Value = 43.40618181228638
$ws.Range("J7").Value = 42.04231452941895
$ws.Range("J8").Value = 40.21286225318909
$ws.Range("J9").Value = 30.53034400939941
$ws.Range("B10").Value = 2073
$ws.Range("E10").Value = 81
$ws.Range("H10").Value = 96.09073359073359
$ws.Range("I10").Value = 0.04306459689534301
$ws.Range("J10").Value = 50.57822561264038
$ws.Range("B11").Value = 2874
$ws.Range("D11").Value = 2873
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 99.9652052887961
$ws.Range("I11").Value = 0.0003478260869565218
$ws.Range("J11").Value = 31.41887331008911
$ws.Range("J12").Value = 41.96691012382507
$ws.Range("J13").Value = 39.74461460113525
$ws.Range("J14").Value = 49.02508592605591
$ws.Range("J15").Value = 36.06474256515503
$ws.Range("J16").Value = 37.9193913936615
$ws.Range("J17").Value = 39.61196422576904
$ws.Range("B18").Value = 2539
$ws.Range("E18").Value = 6
$ws.Range("H18").Value = 99.76359338061465
$ws.Range("I18").Value = 0.003941663381947182
$ws.Range("J18").Value = 50.16433310508728
$ws.Range("J19").Value = 46.6275806427002
$ws.Range("B20").Value = 2940
$ws.Range("D20").Value = 2939
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 99.89802855200544
$ws.Range("I20").Value = 0.001019367991845056
$ws.Range("J20").Value = 42.57859110832214
$ws.Range("B21").Value = 1955
$ws.Range("E21").Value = 5
$ws.Range("H21").Value = 99.74411463664278
$ws.Range("I21").Value = 0.002564102564102564
$ws.Range("J21").Value = 41.69644474983215
$ws.Range("J22").Value = 49.18008327484131
$ws.Range("J23").Value = 43.51635003089905
$ws.Range("J24").Value = 33.82432150840759
$ws.Range("J25").Value = 44.04618239402771
$ws.Range("J26").Value = 42.39572930335999
$ws.Range("B27").Value = 2587
$ws.Range("E27").Value = 3
$ws.Range("H27").Value = 99.88399071925754
$ws.Range("I27").Value = 0.001160990712074303
$ws.Range("J27").Value = 42.55039668083191
$ws.Range("J28").Value = 36.57597899436951
$ws.Range("J29").Value = 42.62749004364014
$ws.Range("J30").Value = 40.26213145256042
$ws.Range("J31").Value = 43.9177393913269
$ws.Range("J32").Value = 36.42376589775085
$ws.Range("B33").Value = 1448
$ws.Range("E33").Value = 2
$ws.Range("H33").Value = 99.86178299930891
$ws.Range("I33").Value = 0.005509641873278237
$ws.Range("J33").Value = 50.54662680625916
$ws.Range("J34").Value = 39.48422980308533
$ws.Range("B35").Value = 3500
$ws.Range("D35").Value = 3497
$ws.Range("F35").Value = 22
$ws.Range("G35").Value = 99.37482239272521
$ws.Range("H35").Value = 99.94284081166047
$ws.Range("I35").Value = 0.006818181818181818
$ws.Range("J35").Value = 32.82644605636597
$ws.Range("J36").Value = 42.92876696586609
$ws.Range("J37").Value = 46.87474322319031
$ws.Range("J38").Value = 33.8890917301178
$ws.Range("B39").Value = 2070
$ws.Range("D39").Value = 2069
$ws.Range("F39").Value = 187
$ws.Range("G39").Value = 91.71099290780141
$ws.Range("I39").Value = 0.08285334514842711
$ws.Range("J39").Value = 34.50550580024719
$ws.Range("B40").Value = 2279
$ws.Range("D40").Value = 2248
$ws.Range("F40").Value = 102
$ws.Range("G40").Value = 95.65957446808511
$ws.Range("H40").Value = 98.68305531167691
$ws.Range("I40").Value = 0.05614632071458953
$ws.Range("J40").Value = 39.72537088394165
$ws.Range("J41").Value = 37.52414703369141
$ws.Range("J42").Value = 39.75350856781006
$ws.Range("J43").Value = 44.16143321990967
$ws.Range("J44").Value = 41.80450344085693
$ws.Range("J45").Value = 37.31579780578613
$ws.Range("J46").Value = 47.26611304283142
$ws.Range("J47").Value = 40.74874758720398
$ws.Range("B48").Value = 2299
$ws.Range("E48").Value = 282
$ws.Range("H48").Value = 87.72845953002611
$ws.Range("I48").Value = 0.1398116013882003
$ws.Range("J48").Value = 45.77669382095337
$ws.Range("J49").Value = 45.7607364654541
$ws.Range("J50").Value = 47.01032614707947
$ws.Range("J51").Value = 40.75990319252014
$ws.Range("B52").Value = 2250
$ws.Range("D52").Value = 2249
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 99.95555555555555
$ws.Range("I52").Value = 0.000444247001332741
$ws.Range("J52").Value = 43.28269004821777
$ws.Range("J53").Value = 45.92657160758972
$ws.Range("J54").Value = 32.86650705337524
$ws.Range("J55").Value = 44.45862507820129
$ws.Range("J56").Value = 36.32990980148315
$ws.Range("J57").Value = 45.01442265510559
$ws.Range("J58").Value = 35.28443479537964
$ws.Range("J59").Value = 34.69120335578918
$ws.Range("J60").Value = 48.14723420143127
$ws.Range("J61").Value = 35.66426277160645
$ws.Range("J62").Value = 36.78492569923401
$ws.Range("B63").Value = 2438
$ws.Range("D63").Value = 2429
$ws.Range("F63").Value = 7
$ws.Range("G63").Value = 99.71264367816092
$ws.Range("H63").Value = 99.6717275338531
$ws.Range("I63").Value = 0.006155108740254411
$ws.Range("J63").Value = 42.70380425453186
$ws.Range("J64").Value = 42.12425804138184
$ws.Range("J65").Value = 38.60061264038086
$ws.Range("J66").Value = 42.69251680374146
$ws.Range("B67").Value = 2808
$ws.Range("D67").Value = 2803
$ws.Range("F67").Value = 93
$ws.Range("G67").Value = 96.78867403314918
$ws.Range("H67").Value = 99.85749910936943
$ws.Range("I67").Value = 0.03348291335864688
$ws.Range("J67").Value = 43.98487067222595
$ws.Range("J68").Value = 44.83152866363525
$ws.Range("J69").Value = 47.04738140106201
$ws.Range("J70").Value = 35.28760194778442
$ws.Range("B71").Value = 1530
$ws.Range("D71").Value = 1529
$ws.Range("F71").Value = 37
$ws.Range("G71").Value = 97.63729246487867
$ws.Range("I71").Value = 0.02361199744735163
$ws.Range("J71").Value = 37.71228075027466
$ws.Range("J72").Value = 43.83773565292358
$ws.Range("B73").Value = 1177
$ws.Range("E73").Value = 48
$ws.Range("H73").Value = 95.91836734693878
$ws.Range("I73").Value = 0.0542432195975503
$ws.Range("J73").Value = 39.01833724975586
$ws.Range("B74").Value = 1327
$ws.Range("D74").Value = 1318
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 37
$ws.Range("G74").Value = 97.26937269372694
$ws.Range("H74").Value = 99.39668174962293
$ws.Range("I74").Value = 0.0331858407079646
$ws.Range("J74").Value = 42.16090941429138
$ws.Range("J75").Value = 43.33949518203735
$ws.Range("B76").Value = 1971
$ws.Range("D76").Value = 1970
$ws.Range("F76").Value = 54
$ws.Range("G76").Value = 97.33201581027669
$ws.Range("I76").Value = 0.02666666666666667
$ws.Range("J76").Value = 37.36109828948975
$ws.Range("J77").Value = 40.39943909645081
